# Commit: "updated till excel part"
# Rewrites the runs/balls/fours/sixes figures for several Shikhar Dhawan (Delhi Capitals)
# innings rows (the per-match stats got reshuffled/corrected).
#
# The source sheet stores these numbers as text (t="str" in the XML, and Excel flags them
# with the "Number Stored as Text" warning) even though the cells use the default General
# format. To keep writing them back as text (instead of letting Excel auto-convert the
# numeric-looking string into a real number), each target cell is switched to the Text ("@")
# number format immediately before its value is assigned.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "32"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "26"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "13"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "0"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "106"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "61"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "12"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "3"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "54"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "41"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "6"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "0"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "35"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "1"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "34"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "31"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "0"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "101"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "14"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "1"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "57"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "6"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "2"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "6"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1"
